# Auto-generated edit script applying the commit diff
# "reworked dur calc, lots of other changes"
$wb = $excel.ActiveWorkbook

# --- site_metrics ---
$ws = $wb.Worksheets.Item("site_metrics")
$ws.Range("O2").Value = 0.02549195762010053
$ws.Range("N5").Value = 8.207563025210083
$ws.Range("O7").Value = 0.004813189684878138
$ws.Range("AK13").Value = $true
$ws.Range("AK14").Value = $true
$ws.Range("O18").Value = 0.009845425644055645
$ws.Range("O30").Value = 0.07811873031498129
$ws.Range("N32").Value = 3.103713151927437
$ws.Range("O33").Value = 0.1402879772438876
$ws.Range("O34").Value = 0.01860710227192392
$ws.Range("AK36").Value = $true
$ws.Range("AK37").Value = $true
$ws.Range("O39").Value = 0.412524850340755
$ws.Range("AK40").Value = $true
$ws.Range("AK41").Value = $true
$ws.Range("AK43").Value = $true
$ws.Range("AK44").Value = $true
$ws.Range("AK52").Value = $true
$ws.Range("O56").Value = 0.08052547057925852
$ws.Range("O57").Value = 0.08906252763741357
$ws.Range("AK60").Value = $true
$ws.Range("N61").Value = 3.469658995974786
$ws.Range("N65").Value = 3.529120879120879
$ws.Range("O67").Value = 0.01099562600072403
$ws.Range("AK69").Value = $true
$ws.Range("N70").Value = 10.93350694444444
$ws.Range("O70").Value = 0.008225483830074253
$ws.Range("Q70").Value = 2.6875
$ws.Range("O71").Value = 0.003885208861282886
$ws.Range("AK73").Value = $true
$ws.Range("N75").Value = 3.402378983024145
$ws.Range("O75").Value = 0.001079330339997736
$ws.Range("N78").Value = 3.49597162097162
$ws.Range("O78").Value = 0.006176649759310394
$ws.Range("O80").Value = 0.004193311531737252

# --- mk_duration ---
$ws = $wb.Worksheets.Item("mk_duration")
$ws.Range("K5").Value = "no trend"
$ws.Range("L5").Value = $false
$ws.Range("M5").Value = 0.09344799392558367
$ws.Range("N5").Value = -1.677483394552118
$ws.Range("O5").Value = -0.2142857142857143
$ws.Range("P5").Value = -87
$ws.Range("Q5").Value = 2628.333333333333
$ws.Range("R5").Value = -0.1091269841269841
$ws.Range("S5").Value = 3.527777777777778
$ws.Range("K26").Value = "no trend"
$ws.Range("L26").Value = $false
$ws.Range("M26").Value = 0.1099636292968564
$ws.Range("N26").Value = 1.598356637186201
$ws.Range("O26").Value = 0.196969696969697
$ws.Range("P26").Value = 104
$ws.Range("Q26").Value = 4152.666666666667
$ws.Range("R26").Value = 0.007905982905982906
$ws.Range("S26").Value = 1.54017094017094
$ws.Range("M31").Value = 0.9523825905788506
$ws.Range("N31").Value = -0.0597150429446484
$ws.Range("O31").Value = -0.006205673758865249
$ws.Range("P31").Value = -7
$ws.Range("Q31").Value = 10095.66666666667
$ws.Range("K32").Value = "no trend"
$ws.Range("L32").Value = $false
$ws.Range("M32").Value = 0.9401083689080136
$ws.Range("N32").Value = 0.07513365721922702
$ws.Range("O32").Value = 0.01231527093596059
$ws.Range("P32").Value = 5
$ws.Range("Q32").Value = 2834.333333333333
$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 2.714285714285714
$ws.Range("K39").Value = "no trend"
$ws.Range("L39").Value = $false
$ws.Range("M39").Value = 0.6157503017697474
$ws.Range("N39").Value = 0.5018823210851089
$ws.Range("O39").Value = 0.05666666666666666
$ws.Range("P39").Value = 17
$ws.Range("Q39").Value = 1016.333333333333
$ws.Range("M61").Value = 0.7214746342177878
$ws.Range("N61").Value = 0.3564886717836127
$ws.Range("O61").Value = 0.03875968992248062
$ws.Range("P61").Value = 35
$ws.Range("Q61").Value = 9096.333333333334
$ws.Range("R61").Value = 0.007575757575757569
$ws.Range("S61").Value = 2.440909090909091
$ws.Range("M65").Value = 0.8123923688472134
$ws.Range("N65").Value = -0.2373408560833091
$ws.Range("O65").Value = -0.03439153439153439
$ws.Range("P65").Value = -13
$ws.Range("Q65").Value = 2556.333333333333
$ws.Range("R65").Value = -0.03003968253968263
$ws.Range("S65").Value = 3.755535714285716
$ws.Range("M66").Value = 0.3660827184679289
$ws.Range("N66").Value = -0.9038353419323745
$ws.Range("O66").Value = -0.1051051051051051
$ws.Range("P66").Value = -70
$ws.Range("Q66").Value = 5828
$ws.Range("R66").Value = -0.03858560794044666
$ws.Range("S66").Value = 3.617617866004963
$ws.Range("D70").Value = 0.6339498028447381
$ws.Range("E70").Value = -0.476174867647586
$ws.Range("F70").Value = -0.06236559139784946
$ws.Range("G70").Value = -29
$ws.Range("I70").Value = -0.05128205128205127
$ws.Range("J70").Value = 8.435897435897436
$ws.Range("M70").Value = 0.2738534999027467
$ws.Range("N70").Value = -1.094231407679569
$ws.Range("O70").Value = -0.1044897959183674
$ws.Range("P70").Value = -128
$ws.Range("M71").Value = 0.4888249449426469
$ws.Range("N71").Value = -0.6921789680474556
$ws.Range("O71").Value = -0.07149758454106281
$ws.Range("P71").Value = -74
$ws.Range("Q71").Value = 11122.66666666667
$ws.Range("R71").Value = -0.01388888888888889
$ws.Range("S71").Value = 2.3125
$ws.Range("M75").Value = 0.5748247878461354
$ws.Range("N75").Value = -0.5609600264316629
$ws.Range("O75").Value = -0.06666666666666667
$ws.Range("P75").Value = -42
$ws.Range("Q75").Value = 5342
$ws.Range("R75").Value = -0.01044383270549599
$ws.Range("S75").Value = 2.18276707234618
$ws.Range("K76").Value = "no trend"
$ws.Range("L76").Value = $false
$ws.Range("M76").Value = 0.5358598692967174
$ws.Range("N76").Value = 0.619085751600795
$ws.Range("O76").Value = 0.06829268292682927
$ws.Range("P76").Value = 56
$ws.Range("Q76").Value = 7892.666666666667
$ws.Range("R76").Value = 0.008130215649012631
$ws.Range("S76").Value = 1.837395687019747
$ws.Range("M78").Value = 0.8869230368864609
$ws.Range("N78").Value = 0.142198727497721
$ws.Range("O78").Value = 0.01545893719806763
$ws.Range("P78").Value = 16
$ws.Range("Q78").Value = 11127.33333333333
$ws.Range("S78").Value = 3.171428571428572
$ws.Range("K80").Value = "no trend"
$ws.Range("L80").Value = $false
$ws.Range("M80").Value = 0.3069912409644038
$ws.Range("N80").Value = 1.021555684791887
$ws.Range("O80").Value = 0.1333333333333333
$ws.Range("P80").Value = 58
$ws.Range("Q80").Value = 3113.333333333333
$ws.Range("R80").Value = 0.02500000000000002
$ws.Range("S80").Value = 2.6375
$ws.Range("M82").Value = 0.858673521284679
$ws.Range("N82").Value = 0.1780629720586193
$ws.Range("O82").Value = 0.0231729055258467
$ws.Range("P82").Value = 13
$ws.Range("Q82").Value = 4541.666666666667
$ws.Range("R82").Value = 0.003654970760233928
$ws.Range("S82").Value = 2.38969298245614

# --- mk_intra_annual ---
$ws = $wb.Worksheets.Item("mk_intra_annual")
$ws.Range("K5").Value = "no trend"
$ws.Range("L5").Value = $false
$ws.Range("M5").Value = 0.410746198342077
$ws.Range("N5").Value = -0.8225811910567383
$ws.Range("O5").Value = -0.1059113300492611
$ws.Range("P5").Value = -43
$ws.Range("Q5").Value = 2607
$ws.Range("S5").Value = 1
$ws.Range("M26").Value = 0.3249311654768372
$ws.Range("N26").Value = -0.984375
$ws.Range("O26").Value = -0.1212121212121212
$ws.Range("P26").Value = -64
$ws.Range("Q26").Value = 4096
$ws.Range("R26").Value = -0.08514492753623187
$ws.Range("S26").Value = 12.36231884057971
$ws.Range("M31").Value = 0.8101656932688071
$ws.Range("N31").Value = -0.2402122813141875
$ws.Range("O31").Value = -0.02216312056737589
$ws.Range("P31").Value = -25
$ws.Range("Q31").Value = 9982.333333333334
$ws.Range("K32").Value = "no trend"
$ws.Range("L32").Value = $false
$ws.Range("M32").Value = 0.9849455743562765
$ws.Range("N32").Value = -0.01886904412032354
$ws.Range("O32").Value = -0.004926108374384237
$ws.Range("P32").Value = -2
$ws.Range("Q32").Value = 2808.666666666667
$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 6
$ws.Range("K39").Value = "no trend"
$ws.Range("L39").Value = $false
$ws.Range("M39").Value = 0.6128224102281419
$ws.Range("N39").Value = 0.5060487741189207
$ws.Range("O39").Value = 0.05666666666666666
$ws.Range("P39").Value = 17
$ws.Range("Q39").Value = 999.6666666666666
$ws.Range("M61").Value = 0.3002074359872218
$ws.Range("N61").Value = -1.035988653519739
$ws.Range("O61").Value = -0.1096345514950166
$ws.Range("P61").Value = -99
$ws.Range("Q61").Value = 8948.333333333334
$ws.Range("R61").Value = -0.03571428571428571
$ws.Range("S61").Value = 5.75
$ws.Range("M65").Value = 0.6755034011772314
$ws.Range("N65").Value = -0.4186069613366101
$ws.Range("O65").Value = -0.0582010582010582
$ws.Range("P65").Value = -22
$ws.Range("Q65").Value = 2516.666666666667
$ws.Range("S65").Value = 4.5
$ws.Range("M66").Value = 0.6251697703125925
$ws.Range("N66").Value = -0.4885366530433578
$ws.Range("O66").Value = -0.05705705705705705
$ws.Range("P66").Value = -38
$ws.Range("Q66").Value = 5736
$ws.Range("S66").Value = 4
$ws.Range("M71").Value = 0.03614983237685898
$ws.Range("N71").Value = 2.095238095238095
$ws.Range("O71").Value = 0.2135265700483092
$ws.Range("P71").Value = 221
$ws.Range("Q71").Value = 11025
$ws.Range("R71").Value = 0.1
$ws.Range("S71").Value = 3.75
$ws.Range("M75").Value = 0.923614755097441
$ws.Range("N75").Value = -0.09588141518882391
$ws.Range("O75").Value = -0.0126984126984127
$ws.Range("P75").Value = -8
$ws.Range("Q75").Value = 5330
$ws.Range("S75").Value = 5
$ws.Range("K76").Value = "no trend"
$ws.Range("L76").Value = $false
$ws.Range("M76").Value = 0.4041663847484855
$ws.Range("N76").Value = 0.8342033836521727
$ws.Range("O76").Value = 0.09146341463414634
$ws.Range("P76").Value = 75
$ws.Range("Q76").Value = 7869
$ws.Range("R76").Value = 0.06559139784946236
$ws.Range("S76").Value = 5.688172043010753
$ws.Range("M78").Value = 0.4913828936852183
$ws.Range("N78").Value = -0.6881109845448636
$ws.Range("O78").Value = -0.07053140096618357
$ws.Range("P78").Value = -73
$ws.Range("Q78").Value = 10948.33333333333
$ws.Range("M80").Value = 0.3146467998598692
$ws.Range("N80").Value = -1.005519426749776
$ws.Range("O80").Value = -0.1310344827586207
$ws.Range("P80").Value = -57
$ws.Range("Q80").Value = 3101.666666666667
$ws.Range("R80").Value = -0.08333333333333333
$ws.Range("S80").Value = 6.708333333333333
$ws.Range("M82").Value = 0.730547120402357
$ws.Range("N82").Value = -0.3443978302354009
$ws.Range("O82").Value = -0.0427807486631016
$ws.Range("P82").Value = -24
$ws.Range("Q82").Value = 4460
$ws.Range("S82").Value = 5

